$d = $word.ActiveDocument

# 1. Discount Rate: 3.0% -> Discount Rate: 3%
$d.Content.Find.Execute("Discount Rate: 3.0%", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Discount Rate: 3%", 2)

# 2. Remove trailing comma after "...internalized into the design process.,"
$d.Content.Find.Execute("internalized into the design process.,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "internalized into the design process.", 2)

# 3. Remove trailing comma after "...won't attract more motorists).," (occurs twice: GHG + Water Pollution)
$d.Content.Find.Execute("more motorists).,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "more motorists).", 2)

# 4. Remove trailing comma after "...sides of the river.,"
$d.Content.Find.Execute("between communities on the two sides of the river.,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "between communities on the two sides of the river.", 2)
